$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5: the financial-documents table switches from the deck's custom
#    "Table_0" table style to the built-in "No Style, Table Grid" style.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{B80DB3D7-A3CD-4068-A27E-88F87EC7936F}")

# ---------------------------------------------------------------------------
# 2) The deck's theme colour palette is switched from the custom "Integral"
#    (Red Violet) palette to the standard Office palette. dk1/lt1 (black and
#    white) are already identical between the two palettes, so only the
#    remaining ten slots need to change.
#    VBA's RGB() packs (r,g,b) into a single BGR-ordered OLE colour value;
#    recreate that helper since it isn't a built-in here.
# ---------------------------------------------------------------------------
function RGB($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$colors = $p.Slides.Range().ThemeColorScheme
$colors.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2       -> 44546A
$colors.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2       -> E7E6E6
$colors.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1   -> 5B9BD5
$colors.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2   -> ED7D31
$colors.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3   -> A5A5A5
$colors.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4   -> FFC000
$colors.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5   -> 4472C4
$colors.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6   -> 70AD47
$colors.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hlink     -> 0563C1
$colors.Colors(12).RGB = RGB 0x95 0x4F 0x72   # folHlink  -> 954F72
